# Fill in most topics
# This script fills in the missing "topic" values (column H) on the
# "outline" worksheet for a large block of rows (6-143), matching the
# values used for similar/adjacent rows in the same subsection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("outline")

# Chapter I: Justice and Security -------------------------------------------------
$ws.Range("H6:H9").Value    = "Corruption Perceptions"
$ws.Range("H10:H18").Value  = "Justice System Evaluation"
$ws.Range("H28").Value      = "Problem Resolution"
$ws.Range("H31:H32").Value  = "Security"
$ws.Range("H33:H39").Value  = "Security Violence"
$ws.Range("H40:H47").Value  = "Law Enforcement Performance"
$ws.Range("H48:H52").Value  = "Criminal Justice Performance"
$ws.Range("H53:H55").Value  = "Law Enforcement Performance"
$ws.Range("H56:H59").Value  = "Criminal Justice Performance"
$ws.Range("H60:H74").Value  = "Perceptions on Authoritarian Behavior"
$ws.Range("H75:H76").Value  = "Justice System Evaluation"
$ws.Range("H77:H79").Value  = "Perceptions on Authoritarian Behavior"

# Civic participation ----------------------------------------------------------
$ws.Range("H80:H81").Value   = "Civic Participation A"
$ws.Range("H82:H85").Value   = "Civic Participation B"
$ws.Range("H86:H88").Value   = "Civic Participation A"
$ws.Range("H93:H94").Value   = "Civic Participation A"
$ws.Range("H95:H97").Value   = "Civic Participation A Civic Participation B"
$ws.Range("H98:H99").Value   = "Civic Participation A"
$ws.Range("H100:H101").Value = "Civic Participation B"

# Discrimination ----------------------------------------------------------------
$ws.Range("H102:H112").Value = "Discrimination"

# Corruption ---------------------------------------------------------------------
$ws.Range("H114").Value      = "Corruption Change"
$ws.Range("H115:H119").Value = "Opinions regarding Corruption"
$ws.Range("H120:H130").Value = "Corruption"
$ws.Range("H131:H135").Value = "Bribe Victimization"

# Access to information -----------------------------------------------------------
$ws.Range("H136:H139").Value = "Information Provision"
$ws.Range("H140:H143").Value = "Information Requests"

# Restore the view state (best effort: last active selection).
$ws.Range("H140").Select()
